# Update betting odds values per the 2025-04-17 FlashScore data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 3.5
$ws.Range("G3").Value = 1.8
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 7
$ws.Range("N3").Value = 2.4
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 1.53
$ws.Range("Q3").Value = 2.38
$ws.Range("Y3").Value = 41
$ws.Range("Z3").Value = 7
$ws.Range("AA3").Value = 6.5
$ws.Range("J5").Value = 1.07
$ws.Range("L5").Value = 1.36
$ws.Range("J6").Value = 1.06
$ws.Range("K6").Value = 10
$ws.Range("L6").Value = 1.33
$ws.Range("G7").Value = 2.5
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 1.1
$ws.Range("L7").Value = 1.44
$ws.Range("M7").Value = 2.63
$ws.Range("N7").Value = 2.5
$ws.Range("O7").Value = 1.5
$ws.Range("P7").Value = 1.57
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 2.1
$ws.Range("S7").Value = 1.67
$ws.Range("W7").Value = 23
$ws.Range("AB7").Value = 19
$ws.Range("G8").Value = 1.67
$ws.Range("J8").Value = 1.07
$ws.Range("K8").Value = 9
$ws.Range("L8").Value = 1.36
$ws.Range("M8").Value = 3
$ws.Range("N8").Value = 2.15
$ws.Range("O8").Value = 1.67
$ws.Range("AF8").Value = 29
$ws.Range("H9").Value = 4.75
$ws.Range("I9").Value = 9
$ws.Range("J9").Value = 1.04
$ws.Range("K9").Value = 12
$ws.Range("L9").Value = 1.22
$ws.Range("M9").Value = 4
$ws.Range("AA9").Value = 9.5
$ws.Range("J10").Value = 1.08
$ws.Range("K10").Value = 8
$ws.Range("N10").Value = 2.4
$ws.Range("O10").Value = 1.53
$ws.Range("G11").Value = 1.91
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 4.2
$ws.Range("J11").Value = 1.1
$ws.Range("K11").Value = 7
$ws.Range("N11").Value = 2.4
$ws.Range("O11").Value = 1.53
$ws.Range("R11").Value = 2.2
$ws.Range("S11").Value = 1.62
$ws.Range("T11").Value = 5.5
$ws.Range("V11").Value = 9.5
$ws.Range("Y11").Value = 41
$ws.Range("Z11").Value = 7
$ws.Range("AB11").Value = 21
$ws.Range("AC11").Value = 81
$ws.Range("AE11").Value = 9
$ws.Range("AF11").Value = 19
$ws.Range("AH11").Value = 41
$ws.Range("AA12").Value = 7.7
$ws.Range("G13").Value = 2.15
$ws.Range("I13").Value = 3.3
$ws.Range("T13").Value = 8
$ws.Range("U13").Value = 11
$ws.Range("W13").Value = 21
$ws.Range("AH13").Value = 34
$ws.Range("AI13").Value = 26
$ws.Range("G14").Value = 2.12
$ws.Range("H14").Value = 3.4
$ws.Range("I14").Value = 3.1
$ws.Range("T14").Value = 7.6
$ws.Range("V14").Value = 8.75
$ws.Range("W14").Value = 19.5
$ws.Range("X14").Value = 17
$ws.Range("Y14").Value = 28
$ws.Range("AA14").Value = 6.5
$ws.Range("AB14").Value = 14.5
$ws.Range("AE14").Value = 9.5
$ws.Range("AG14").Value = 11.25
$ws.Range("AI14").Value = 28
$ws.Range("AJ14").Value = 37
$ws.Range("G15").Value = 3.95
$ws.Range("H15").Value = 3.7
$ws.Range("I15").Value = 1.75
$ws.Range("N15").Value = 1.57
$ws.Range("O15").Value = 2.12
$ws.Range("T15").Value = 14
$ws.Range("U15").Value = 25
$ws.Range("V15").Value = 13
$ws.Range("W15").Value = 60
$ws.Range("X15").Value = 35
$ws.Range("Y15").Value = 35
$ws.Range("Z15").Value = 14
$ws.Range("AA15").Value = 7.5
$ws.Range("AD15").Value = 300
$ws.Range("AE15").Value = 9.25
$ws.Range("AF15").Value = 10
$ws.Range("AG15").Value = 8.25
$ws.Range("AH15").Value = 15.5
$ws.Range("AI15").Value = 12.5
$ws.Range("AJ15").Value = 19.5
$ws.Range("L16").Value = 1.2
$ws.Range("M16").Value = 4.33
$ws.Range("N16").Value = 1.67
$ws.Range("O16").Value = 2.15
$ws.Range("G17").Value = 1.33
$ws.Range("H17").Value = 5.25
$ws.Range("I17").Value = 9
$ws.Range("L17").Value = 1.17
$ws.Range("M17").Value = 5
$ws.Range("N17").Value = 1.57
$ws.Range("O17").Value = 2.35
$ws.Range("U17").Value = 7
$ws.Range("AC17").Value = 51
$ws.Range("AD17").Value = 301
$ws.Range("AE17").Value = 21
$ws.Range("AG17").Value = 23
$ws.Range("G18").Value = 2.4
$ws.Range("I18").Value = 3
$ws.Range("L18").Value = 1.33
$ws.Range("M18").Value = 3.25
$ws.Range("U18").Value = 12
$ws.Range("V18").Value = 10
$ws.Range("W18").Value = 23
$ws.Range("X18").Value = 21
$ws.Range("AE18").Value = 9.5
$ws.Range("AG18").Value = 11
$ws.Range("AH18").Value = 29
$ws.Range("AI18").Value = 23
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 3.4
$ws.Range("I19").Value = 3.4
$ws.Range("J19").Value = 1.02
$ws.Range("K19").Value = 11
$ws.Range("U19").Value = 11
$ws.Range("W19").Value = 19
$ws.Range("X19").Value = 17
$ws.Range("Z19").Value = 11
$ws.Range("AA19").Value = 6.5
$ws.Range("AH19").Value = 34
$ws.Range("AI19").Value = 26
$ws.Range("K20").Value = 17
$ws.Range("L20").Value = 1.06
$ws.Range("M20").Value = 8
$ws.Range("N20").Value = 1.25
$ws.Range("O20").Value = 3.75
$ws.Range("R20").Value = 1.73
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 15
$ws.Range("AA20").Value = 17
$ws.Range("G21").Value = 1.53
$ws.Range("H21").Value = 4.1
$ws.Range("I21").Value = 5.5
$ws.Range("N21").Value = 1.6
$ws.Range("O21").Value = 2.3
$ws.Range("U21").Value = 8.5
$ws.Range("W21").Value = 12
$ws.Range("Y21").Value = 21
$ws.Range("AA21").Value = 8
$ws.Range("AF21").Value = 29
$ws.Range("AG21").Value = 17
$ws.Range("G22").Value = 3.7
$ws.Range("I22").Value = 1.8
$ws.Range("L22").Value = 1.1
$ws.Range("Y22").Value = 23
$ws.Range("AB22").Value = 11
$ws.Range("AI22").Value = 13
$ws.Range("G23").Value = 2.25
$ws.Range("H23").Value = 3.05
$ws.Range("I23").Value = 3.15
$ws.Range("R23").Value = 1.98
$ws.Range("S23").Value = 1.65
$ws.Range("T23").Value = 6.1
$ws.Range("V23").Value = 9.5
$ws.Range("X23").Value = 22
$ws.Range("Y23").Value = 40
$ws.Range("Z23").Value = 7.1
$ws.Range("AA23").Value = 6
$ws.Range("AE23").Value = 7.8
$ws.Range("AF23").Value = 15
$ws.Range("AG23").Value = 11.75
$ws.Range("AH23").Value = 40
$ws.Range("AI23").Value = 32

Write-Output "Updated 193 cells"
